$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix capitalization of the two image-description cells (first letter of each
# sentence was lower-case in the original upload).
$ws.Range("A45").Value = "Man in a black shirt with short hair is installing or adjusting a metal panel on the underside of an aircraft. Frontal and close-up photo of the man."
$ws.Range("A48").Value = "White-haired man stepped out, holding a black-framed glass door in his hand. In front of him was a multicolored flag, including red, yellow, and green. Several microphones were raised near the man. Among them was a green microphone. A close-up of the man's face."

# Select the whole used range and scroll the window down so row 68 is at the
# top, matching the author's final view state when the file was saved.
$ws.Range("A1:A95").Select()
$excel.ActiveWindow.ScrollRow = 68
